$d = $word.ActiveDocument

# 1. Insert the three new log lines as new paragraphs after the current last
#    paragraph's text (this lands the new text before the existing _GoBack
#    bookmark, which stays attached to "...Filetered DS/Flow clicked").
$lastPara = $d.Paragraphs.Last
$lastPara.Range.InsertAfter("`rD_Monitor-TCNo_001_SUCCESS: Running and Completed tab is displayed in monitor screen`rD_Monitor-TCNo_002_WARN: UiElement_FilteredResult.IsVisibleFalse`rD_Monitor-TCNo_004_SUCCESS: Ui checked for Completed Tab")

# 2. Relocate the _GoBack bookmark so that it sits, zero-width, right after
#    the text of the new last paragraph (matching the target XML shape).
#    Directly adding a collapsed bookmark exactly at "end of paragraph text"
#    is unreliable, so a temporary marker character is used to get a safe,
#    non-boundary insertion point, then removed again.
$newLastPara = $d.Paragraphs.Last
$insertPos = $newLastPara.Range.End - 1

$marker = $d.Range($insertPos, $insertPos)
$marker.InsertAfter("~")

$bmRange = $d.Range($insertPos, $insertPos)
$d.Bookmarks.Add("_GoBack", $bmRange)

$markerRange = $d.Range($insertPos, $insertPos + 1)
$markerRange.Delete()
